$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: pre-seed formatting (no explicit style) for the whole new block by copying
# the format of the last existing data row (653), which itself carries no explicit style.
$ws.Range("A653").Copy() | Out-Null
$ws.Range("A654:A673").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C653:J653").Copy() | Out-Null
$ws.Range("C654:J673").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Step 2: write the appended data rows (dates as literal text, matching the existing column A cells)
$ws.Range("A654").NumberFormat = "@"
$ws.Range("A654").Value = "2024-09-02"
$ws.Range("C654").Value = 650.9500122070312
$ws.Range("D654").Value = 1470.050048828125
$ws.Range("E654").Value = 608.5800170898438
$ws.Range("F654").Value = 1316.800048828125
$ws.Range("G654").Value = 811.2000122070312
$ws.Range("H654").Value = 17987.74047851562
$ws.Range("I654").Value = 0
$ws.Range("J654").Value = 186.5558764959666

$ws.Range("A655").NumberFormat = "@"
$ws.Range("A655").Value = "2024-09-03"
$ws.Range("C655").Value = 640.0499877929688
$ws.Range("D655").Value = 1460.75
$ws.Range("E655").Value = 599.9400024414062
$ws.Range("F655").Value = 1341.949951171875
$ws.Range("G655").Value = 822.3499755859375
$ws.Range("H655").Value = 17977.66967773438
$ws.Range("I655").Value = -0.0005598702512568748
$ws.Range("J655").Value = 186.4514294105194

$ws.Range("A656").NumberFormat = "@"
$ws.Range("A656").Value = "2024-09-04"
$ws.Range("C656").Value = 645.5999755859375
$ws.Range("D656").Value = 1475.300048828125
$ws.Range("E656").Value = 609
$ws.Range("F656").Value = 1327.75
$ws.Range("G656").Value = 824.2999877929688
$ws.Range("H656").Value = 18052.54992675781
$ws.Range("I656").Value = 0.004165181047695957
$ws.Range("J656").Value = 187.2280333706159

$ws.Range("A657").NumberFormat = "@"
$ws.Range("A657").Value = "2024-09-05"
$ws.Range("C657").Value = 643.9000244140625
$ws.Range("D657").Value = 1457.699951171875
$ws.Range("E657").Value = 602.1799926757812
$ws.Range("F657").Value = 1254.800048828125
$ws.Range("G657").Value = 835.4000244140625
$ws.Range("H657").Value = 17792.94024658203
$ws.Range("I657").Value = -0.01438077619112318
$ws.Range("J657").Value = 184.5355489260089

$ws.Range("A658").NumberFormat = "@"
$ws.Range("A658").Value = "2024-09-06"
$ws.Range("C658").Value = 665.25
$ws.Range("D658").Value = 1443.449951171875
$ws.Range("E658").Value = 597.2999877929688
$ws.Range("F658").Value = 1256.849975585938
$ws.Range("G658").Value = 832.7000122070312
$ws.Range("H658").Value = 17880.34979248047
$ws.Range("I658").Value = 0.004912597057432855
$ws.Range("J658").Value = 185.4420977206545

$ws.Range("A659").NumberFormat = "@"
$ws.Range("A659").Value = "2024-09-09"
$ws.Range("C659").Value = 676
$ws.Range("D659").Value = 1492.050048828125
$ws.Range("E659").Value = 610.3400268554688
$ws.Range("F659").Value = 1225.25
$ws.Range("G659").Value = 827.5999755859375
$ws.Range("H659").Value = 18025.32012939453
$ws.Range("I659").Value = 0.008107802061849448
$ws.Range("J659").Value = 186.9456255429077

$ws.Range("A660").NumberFormat = "@"
$ws.Range("A660").Value = "2024-09-10"
$ws.Range("C660").Value = 680
$ws.Range("D660").Value = 1503.050048828125
$ws.Range("E660").Value = 608
$ws.Range("F660").Value = 1246
$ws.Range("G660").Value = 824.75
$ws.Range("H660").Value = 18130.15014648438
$ws.Range("I660").Value = 0.005815709032478913
$ws.Range("J660").Value = 188.03284690596

$ws.Range("A661").NumberFormat = "@"
$ws.Range("A661").Value = "2024-09-11"
$ws.Range("C661").Value = 680.4500122070312
$ws.Range("D661").Value = 1499.949951171875
$ws.Range("E661").Value = 627.6599731445312
$ws.Range("F661").Value = 1229
$ws.Range("G661").Value = 814
$ws.Range("H661").Value = 18088.97985839844
$ws.Range("I661").Value = -0.002270818926114677
$ws.Range("J661").Value = 187.6058583584748

$ws.Range("A662").NumberFormat = "@"
$ws.Range("A662").Value = "2024-09-12"
$ws.Range("C662").Value = 686.0999755859375
$ws.Range("D662").Value = 1513.449951171875
$ws.Range("E662").Value = 645.5999755859375
$ws.Range("F662").Value = 1224.849975585938
$ws.Range("G662").Value = 809.7000122070312
$ws.Range("H662").Value = 18193.19958496094
$ws.Range("I662").Value = 0.00576150382046627
$ws.Range("J662").Value = 188.686750228149

$ws.Range("A663").NumberFormat = "@"
$ws.Range("A663").Value = "2024-09-13"
$ws.Range("C663").Value = 681.9500122070312
$ws.Range("D663").Value = 1491.300048828125
$ws.Range("E663").Value = 646.6500244140625
$ws.Range("F663").Value = 1229.300048828125
$ws.Range("G663").Value = 788.0499877929688
$ws.Range("H663").Value = 18027.60040283203
$ws.Range("I663").Value = -0.009102257211853799
$ws.Range("J663").Value = 186.9692748951036

$ws.Range("A664").NumberFormat = "@"
$ws.Range("A664").Value = "2024-09-16"
$ws.Range("C664").Value = 695.2000122070312
$ws.Range("D664").Value = 1456.349975585938
$ws.Range("E664").Value = 621.0499877929688
$ws.Range("F664").Value = 1219.699951171875
$ws.Range("G664").Value = 751.9500122070312
$ws.Range("H664").Value = 17765.49987792969
$ws.Range("I664").Value = -0.01453884704817227
$ws.Range("J664").Value = 184.250957204696

$ws.Range("A665").NumberFormat = "@"
$ws.Range("A665").Value = "2024-09-17"
$ws.Range("C665").Value = 692
$ws.Range("D665").Value = 1459.400024414062
$ws.Range("E665").Value = 649.6500244140625
$ws.Range("F665").Value = 1222.949951171875
$ws.Range("G665").Value = 746.75
$ws.Range("H665").Value = 17827
$ws.Range("I665").Value = 0.003461772676980224
$ws.Range("J665").Value = 184.8887921340547

$ws.Range("A666").NumberFormat = "@"
$ws.Range("A666").Value = "2024-09-18"
$ws.Range("C666").Value = 695.2999877929688
$ws.Range("D666").Value = 1432.150024414062
$ws.Range("E666").Value = 646.7000122070312
$ws.Range("F666").Value = 1224.550048828125
$ws.Range("G666").Value = 744.5999755859375
$ws.Range("H666").Value = 17755.70007324219
$ws.Range("I666").Value = -0.003999547133999692
$ws.Range("J666").Value = 184.1493206953662

$ws.Range("A667").NumberFormat = "@"
$ws.Range("A667").Value = "2024-09-19"
$ws.Range("C667").Value = 697
$ws.Range("D667").Value = 1444.849975585938
$ws.Range("E667").Value = 652.1500244140625
$ws.Range("F667").Value = 1197.849975585938
$ws.Range("G667").Value = 747.2000122070312
$ws.Range("H667").Value = 17752.34997558594
$ws.Range("I667").Value = -0.0001886773060161447
$ws.Range("J667").Value = 184.1145758976327

$ws.Range("A668").NumberFormat = "@"
$ws.Range("A668").Value = "2024-09-20"
$ws.Range("C668").Value = 709
$ws.Range("D668").Value = 1456.599975585938
$ws.Range("E668").Value = 654.4500122070312
$ws.Range("F668").Value = 1206.300048828125
$ws.Range("G668").Value = 747.5499877929688
$ws.Range("H668").Value = 17905.25006103516
$ws.Range("I668").Value = 0.008612949027001824
$ws.Range("J668").Value = 185.7003453549671

$ws.Range("A669").NumberFormat = "@"
$ws.Range("A669").Value = "2024-09-23"
$ws.Range("C669").Value = 702.5
$ws.Range("D669").Value = 1449.300048828125
$ws.Range("E669").Value = 654.0999755859375
$ws.Range("F669").Value = 1190
$ws.Range("G669").Value = 763.75
$ws.Range("H669").Value = 17852.70007324219
$ws.Range("I669").Value = -0.00293489270542646
$ws.Range("J669").Value = 185.1553347659896

$ws.Range("A670").NumberFormat = "@"
$ws.Range("A670").Value = "2024-09-24"
$ws.Range("C670").Value = 705.0999755859375
$ws.Range("D670").Value = 1446.349975585938
$ws.Range("E670").Value = 646.8499755859375
$ws.Range("F670").Value = 1194.699951171875
$ws.Range("G670").Value = 760.9500122070312
$ws.Range("H670").Value = 17843.19958496094
$ws.Range("I670").Value = -0.0005321597429113499
$ws.Range("J670").Value = 185.0568025506419

$ws.Range("A671").NumberFormat = "@"
$ws.Range("A671").Value = "2024-09-25"
$ws.Range("C671").Value = 689.2000122070312
$ws.Range("D671").Value = 1429.550048828125
$ws.Range("E671").Value = 633.2999877929688
$ws.Range("F671").Value = 1175.349975585938
$ws.Range("G671").Value = 742.5499877929688
$ws.Range("H671").Value = 17509.20007324219
$ws.Range("I671").Value = -0.0187185885652627
$ws.Range("J671").Value = 181.5928004024934

$ws.Range("A672").NumberFormat = "@"
$ws.Range("A672").Value = "2024-09-26"
$ws.Range("C672").Value = 693.5999755859375
$ws.Range("D672").Value = 1422.300048828125
$ws.Range("E672").Value = 626.8499755859375
$ws.Range("F672").Value = 1165
$ws.Range("G672").Value = 742.25
$ws.Range("H672").Value = 17466.64990234375
$ws.Range("I672").Value = -0.002430160756656341
$ws.Range("J672").Value = 181.1515007052639

$ws.Range("A673").NumberFormat = "@"
$ws.Range("A673").Value = "2024-09-27"
$ws.Range("C673").Value = 692.4500122070312
$ws.Range("D673").Value = 1388.650024414062
$ws.Range("E673").Value = 608.5499877929688
$ws.Range("F673").Value = 1165.550048828125
$ws.Range("G673").Value = 751.6500244140625
$ws.Range("H673").Value = 17342.00036621094
$ws.Range("I673").Value = -0.007136430673868747
$ws.Range("J673").Value = 179.8587255790135

# Step 3: re-copy the clean format over column A once more so the temporary Text
# NumberFormat ("@") we applied to coerce the dates into literal strings does not
# leave a stray style index on the new cells (matches rows 634-653 which have none).
$ws.Range("A653").Copy() | Out-Null
$ws.Range("A654:A673").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = 0
